$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")
$ws.Range("E5").Value = 153
$ws.Range("F5").Value = 105
$ws.Range("H5").Value = 116
$ws.Range("E6").Value = 47
$ws.Range("E7").Value = 39
$ws.Range("F7").Value = 24
$ws.Range("H7").Value = 28
$ws.Range("E9").Value = 11
$ws.Range("E10").Value = 646
$ws.Range("F10").Value = 348
$ws.Range("H10").Value = 444
$ws.Range("E11").Value = 422
$ws.Range("F11").Value = 230
$ws.Range("H11").Value = 294
$ws.Range("E12").Value = 650
$ws.Range("F12").Value = 378
$ws.Range("H12").Value = 464
$ws.Range("E14").Value = 134
$ws.Range("E15").Value = 189
$ws.Range("E16").Value = 225
$ws.Range("F16").Value = 125
$ws.Range("H16").Value = 173
$ws.Range("E17").Value = 119
$ws.Range("F17").Value = 65
$ws.Range("H17").Value = 89
$ws.Range("E18").Value = 56
$ws.Range("F18").Value = 29
$ws.Range("H18").Value = 46
$ws.Range("F20").Value = 38
$ws.Range("H20").Value = 75
$ws.Range("F21").Value = 88
$ws.Range("H21").Value = 119
$ws.Range("E22").Value = 184
$ws.Range("F22").Value = 103
$ws.Range("H22").Value = 145
$ws.Range("E23").Value = 215
$ws.Range("E24").Value = 246
$ws.Range("F24").Value = 142
$ws.Range("H24").Value = 172
$ws.Range("E25").Value = 306
$ws.Range("F25").Value = 166
$ws.Range("H25").Value = 226
$ws.Range("E26").Value = 178
$ws.Range("F26").Value = 108
$ws.Range("H26").Value = 133
$ws.Range("E27").Value = 365
$ws.Range("F27").Value = 195
$ws.Range("H27").Value = 276
$ws.Range("E28").Value = 219
$ws.Range("F28").Value = 104
$ws.Range("H28").Value = 156
$ws.Range("E29").Value = 185
$ws.Range("E30").Value = 243
$ws.Range("F30").Value = 148
$ws.Range("H30").Value = 200
$ws.Range("E31").Value = 80
$ws.Range("E32").Value = 200
$ws.Range("F32").Value = 127
$ws.Range("H32").Value = 165
$ws.Range("E33").Value = 317
$ws.Range("E34").Value = 241
$ws.Range("F34").Value = 164
$ws.Range("H34").Value = 202
$ws.Range("E35").Value = 172
$ws.Range("F35").Value = 116
$ws.Range("H35").Value = 143
$ws.Range("E36").Value = 87
$ws.Range("F36").Value = 53
$ws.Range("H36").Value = 63
$ws.Range("E37").Value = 183
$ws.Range("E38").Value = 100
$ws.Range("E39").Value = 190
$ws.Range("F39").Value = 99
$ws.Range("H39").Value = 150
$ws.Range("E40").Value = 288
$ws.Range("F40").Value = 143
$ws.Range("H40").Value = 223
$ws.Range("E41").Value = 422
$ws.Range("E42").Value = 431
$ws.Range("F42").Value = 242
$ws.Range("H42").Value = 303
$ws.Range("E44").Value = 345
$ws.Range("F44").Value = 179
$ws.Range("H44").Value = 247
$ws.Range("E45").Value = 170
$ws.Range("F45").Value = 89
$ws.Range("H45").Value = 128
$ws.Range("E46").Value = 369
$ws.Range("F46").Value = 202
$ws.Range("H46").Value = 265
$ws.Range("E47").Value = 514
$ws.Range("F47").Value = 279
$ws.Range("H47").Value = 371
$ws.Range("E48").Value = 250
$ws.Range("F48").Value = 114
$ws.Range("H48").Value = 158
$ws.Range("E49").Value = 321
$ws.Range("E50").Value = 266
$ws.Range("F50").Value = 140
$ws.Range("H50").Value = 211
$ws.Range("E52").Value = 31
$ws.Range("F52").Value = 14
$ws.Range("H52").Value = 22